$d = $word.ActiveDocument

# The document status line reads "Status:  Proposed" -> change to "Released".
# (Two other "Proposed" values appear as historical status entries in the
# revision-history table further down and must stay untouched, so only the
# first match is replaced - wdReplaceOne (1), not wdReplaceAll.)
$d.Content.Find.Execute("Proposed", $false, $true, $false, $false, $false, $true, 1, $false, "Released", 1) | Out-Null

# Word keeps a single hidden "_GoBack" bookmark marking the location of the
# most recent edit. Re-create it split across "Relea|sed" (where the typed
# replacement happened) - this also removes whatever "_GoBack" bookmark
# existed before (the one around "Req_PO1_DGC_CYRS_012...v1.0") and bumps
# any bookmark ids that collide with the freshly minted one.
$rng = $d.Content
$rng.Find.Execute("Released", $false, $true, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$goBackPos = $rng.Start + 5
$goBackRange = $d.Range($goBackPos, $goBackPos)
$d.Bookmarks.Add("_GoBack", $goBackRange) | Out-Null
